$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (D)
$ws.Range("B3").Value = 1.240874295700121
$ws.Range("C3").Value = 0.9676903269375425
$ws.Range("D3").Value = 0.9262671631783286

# Row 5 (D2)
$ws.Range("B5").Value = 6.103355425606415
$ws.Range("C5").Value = 20.0138924463398
$ws.Range("D5").Value = 61.25608461685

# Row 8 (lambda == 1)
$ws.Range("B8").Value = 15.71438568301618
$ws.Range("C8").Value = 126.7731092634314
$ws.Range("D8").Value = 2174.868861313473

# Row 9 (RES_2_1)
$ws.Range("B9").Value = "H0"

# Row 10 (lambda == 1.3)
$ws.Range("B10").Value = 118.3031466402939
$ws.Range("C10").Value = 3936.864099980687
$ws.Range("D10").Value = 31766.96556149196

# Row 13 (mu 1)
$ws.Range("B13").Value = 58
$ws.Range("C13").Value = 642
$ws.Range("D13").Value = 6782

# Row 16 (mu 1.2)
$ws.Range("B16").Value = 81
$ws.Range("C16").Value = 744
$ws.Range("D16").Value = 7535

# Row 20 (lambda: 1)
$ws.Range("B20").Value = 0.04199999999999998
$ws.Range("C20").Value = 0.01689999999999992
$ws.Range("D20").Value = 0.004319999999999991

# Row 24 (lambda: 1.2)
$ws.Range("B24").Value = 0.1000000000000001
$ws.Range("C24").Value = 0.1073000000000001
$ws.Range("D24").Value = 0.09263999999999994
